$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "min load time" strings and new "min load val" numbers for rows 2-11
# (row index in worksheet -> A=index, B=value, C=time string)
$times = @("18:06:00", "18:33:00", "18:49:00", "18:51:00", "18:33:00", "18:45:00", "18:42:00", "18:32:00", "18:39:00", "18:03:00")
$vals  = @(1.02484, 1.02436, 1.02272, 1.02507, 1.02475, 1.02391, 1.02554, 1.02534, 1.02481, 1.02533)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $vals[$i]
    $ws.Cells.Item($row, 3).Value = $times[$i]
}
